$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 32   Number  22"
$ws.Range("C9").Value = "Report Covering the Week  5/26/2025  Through  6/1/2025"

# --- Plain numeric value updates ---
$ws.Range("F15").Value = 2
$ws.Range("M15").Value = 87.5
$ws.Range("N15").Value = 50
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 300
$ws.Range("F16").Value = 22
$ws.Range("H16").Value = 175
$ws.Range("I16").Value = 100
$ws.Range("J16").Value = 91
$ws.Range("K16").Value = 9.890109890109
$ws.Range("L16").Value = 31.578947368421
$ws.Range("M16").Value = 1.010101010101
$ws.Range("N16").Value = -46.808510638297
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 25
$ws.Range("F17").Value = 44
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = 51.724137931034
$ws.Range("I17").Value = 158
$ws.Range("J17").Value = 141
$ws.Range("K17").Value = 12.056737588652
$ws.Range("L17").Value = 29.508196721311
$ws.Range("M17").Value = 105.194805194805
$ws.Range("N17").Value = 71.739130434782
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 2
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -36.363636363636
$ws.Range("I18").Value = 64
$ws.Range("J18").Value = 61
$ws.Range("K18").Value = 4.918032786885
$ws.Range("L18").Value = 3.225806451612
$ws.Range("M18").Value = -49.606299212598
$ws.Range("N18").Value = -82.465753424657
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 53
$ws.Range("G19").Value = 61
$ws.Range("H19").Value = -13.114754098360
$ws.Range("I19").Value = 286
$ws.Range("J19").Value = 274
$ws.Range("K19").Value = 4.379562043795
$ws.Range("L19").Value = 1.418439716312
$ws.Range("M19").Value = 73.333333333333
$ws.Range("N19").Value = 59.776536312849
$ws.Range("C20").Value = 15
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = 87.5
$ws.Range("F20").Value = 47
$ws.Range("G20").Value = 30
$ws.Range("H20").Value = 56.666666666666
$ws.Range("I20").Value = 221
$ws.Range("J20").Value = 162
$ws.Range("K20").Value = 36.419753086419
$ws.Range("L20").Value = -4.329004329004
$ws.Range("M20").Value = 211.267605633803
$ws.Range("N20").Value = -72.512437810945
$ws.Range("C21").Value = 40
$ws.Range("D21").Value = 34
$ws.Range("E21").Value = 17.647058823529
$ws.Range("F21").Value = 175
$ws.Range("G21").Value = 139
$ws.Range("H21").Value = 25.899280575539
$ws.Range("I21").Value = 847
$ws.Range("J21").Value = 740
$ws.Range("K21").Value = 14.459459459459
$ws.Range("L21").Value = 8.589743589743
$ws.Range("M21").Value = 54.844606946983
$ws.Range("N21").Value = -48.385131017672
$ws.Range("L22").Value = -66.666666666666
$ws.Range("C23").Value = 1
$ws.Range("F23").Value = 10
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 400
$ws.Range("I23").Value = 35
$ws.Range("K23").Value = 34.615384615384
$ws.Range("L23").Value = 40
$ws.Range("M23").Value = 118.75
$ws.Range("C24").Value = 37
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 105.555555555556
$ws.Range("F24").Value = 93
$ws.Range("G24").Value = 84
$ws.Range("H24").Value = 10.714285714285
$ws.Range("I24").Value = 596
$ws.Range("J24").Value = 570
$ws.Range("K24").Value = 4.561403508771
$ws.Range("L24").Value = 3.832752613240
$ws.Range("M24").Value = 12.030075187969
$ws.Range("C25").Value = 20
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 122.222222222222
$ws.Range("F25").Value = 46
$ws.Range("G25").Value = 50
$ws.Range("H25").Value = -8
$ws.Range("I25").Value = 321
$ws.Range("J25").Value = 329
$ws.Range("K25").Value = -2.431610942249
$ws.Range("L25").Value = -3.603603603603
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 42
$ws.Range("G26").Value = 35
$ws.Range("H26").Value = 20
$ws.Range("I26").Value = 208
$ws.Range("J26").Value = 223
$ws.Range("K26").Value = -6.726457399103
$ws.Range("L26").Value = -2.803738317757
$ws.Range("M26").Value = 26.060606060606
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 100
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 5
$ws.Range("H28").Value = 25
$ws.Range("I28").Value = 28
$ws.Range("J28").Value = 33
$ws.Range("K28").Value = -15.151515151515
$ws.Range("L28").Value = 64.705882352941
$ws.Range("J29").Value = 3
$ws.Range("K29").Value = 33.333333333333
$ws.Range("N29").Value = -42.857142857142
$ws.Range("J30").Value = 3
$ws.Range("K30").Value = 33.333333333333
$ws.Range("N30").Value = -42.857142857142
$ws.Range("L33").Value = -50

# --- Cells changing type (number <-> shared text "0"/"***.*") ---
# Each: set raw value (apostrophe-prefixed for text) then paste format from a stable donor cell of the right style
$ws.Range("G14").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("H14").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("G15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("H15").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("C22").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D23").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("C28").Value = 2
$ws.Range("I14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D29").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("D30").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
